$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2
$wdFindContinue = 1
$wdReplaceAll = 2

# --- 1. Title-page "Date:" line (body content) ---
$d.Content.Find.Execute(
    "Date: 11:17 06 Nov 2025 ", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "Date: 11:38 06 Nov 2025 ", $wdReplaceAll
)

# --- 2. The two "generated at" timestamps inside the log-file table (body content) ---
$d.Content.Find.Execute(
    "2025-11-06 11:16:12", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "2025-11-06 11:37:45", $wdReplaceAll
)

# --- 3. Footer text repeated in each section's footer ---
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)
    for ($f = 1; $f -le $section.Footers.Count; $f++) {
        $footer = $section.Footers.Item($f)
        if ($footer.Exists) {
            $footer.Range.Find.Execute(
                "Cambridge CTU, Simon Bond - 11:17 06 Nov 2025 - Page ", $true, $false, $false, $false, $false,
                $true, $wdFindContinue, $false, "Cambridge CTU, Simon Bond - 11:38 06 Nov 2025 - Page ", $wdReplaceAll
            )
        }
    }
}

Write-Output "done"
